$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.2989916666666667
$ws.Cells.Item(2, 8).Value = 0.896975
$ws.Cells.Item(2, 9).Value = 0.01120651476222736
$ws.Cells.Item(2, 10).Value = 0.01120651476222736
$ws.Cells.Item(2, 13).Value = 0.110552
$ws.Cells.Item(2, 14).Value = 0.331656
$ws.Cells.Item(2, 15).Value = 0.01126249561724847
$ws.Cells.Item(2, 16).Value = 0.01126249561724847
$ws.Cells.Item(2, 17).Value = 0.03305412673333333
$ws.Cells.Item(2, 18).Value = 0.2974871406
$ws.Cells.Item(2, 19).Value = 0.0001262133233942159
$ws.Cells.Item(2, 20).Value = 0.000126213323394216

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.2989916666666667
$ws.Cells.Item(3, 8).Value = 0.896975
$ws.Cells.Item(3, 9).Value = 0.01120651476222736
$ws.Cells.Item(3, 10).Value = 0.01120651476222736
$ws.Cells.Item(3, 15).Value = 0.9181055646724333
$ws.Cells.Item(3, 16).Value = 0.9181055646724334
$ws.Cells.Item(3, 17).Value = 2.694534028744444
$ws.Cells.Item(3, 18).Value = 24.2508062587
$ws.Cells.Item(3, 19).Value = 0.01028876356378471
$ws.Cells.Item(3, 20).Value = 0.01028876356378471

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.2989916666666667
$ws.Cells.Item(4, 8).Value = 0.896975
$ws.Cells.Item(4, 9).Value = 0.01120651476222736
$ws.Cells.Item(4, 10).Value = 0.01120651476222736
$ws.Cells.Item(4, 13).Value = 0.6933189999999999
$ws.Cells.Item(4, 14).Value = 2.079957
$ws.Cells.Item(4, 15).Value = 0.07063193971031816
$ws.Cells.Item(4, 16).Value = 0.07063193971031817
$ws.Cells.Item(4, 17).Value = 0.2072966033416666
$ws.Cells.Item(4, 18).Value = 1.865669430075
$ws.Cells.Item(4, 19).Value = 0.0007915378750484333
$ws.Cells.Item(4, 20).Value = 0.0007915378750484335

$ws.Cells.Item(5, 9).Value = 0.9329357354307251
$ws.Cells.Item(5, 10).Value = 0.9329357354307249
$ws.Cells.Item(5, 13).Value = 0.110552
$ws.Cells.Item(5, 14).Value = 0.331656
$ws.Cells.Item(5, 15).Value = 0.01126249561724847
$ws.Cells.Item(5, 16).Value = 0.01126249561724847
$ws.Cells.Item(5, 17).Value = 2.751736528909333
$ws.Cells.Item(5, 18).Value = 24.765628760184
$ws.Cells.Item(5, 19).Value = 0.01050718463146302
$ws.Cells.Item(5, 20).Value = 0.01050718463146302

$ws.Cells.Item(6, 9).Value = 0.9329357354307251
$ws.Cells.Item(6, 10).Value = 0.9329357354307249
$ws.Cells.Item(6, 15).Value = 0.9181055646724333
$ws.Cells.Item(6, 16).Value = 0.9181055646724334
$ws.Cells.Item(6, 19).Value = 0.8565334901807177
$ws.Cells.Item(6, 20).Value = 0.8565334901807177

$ws.Cells.Item(7, 9).Value = 0.9329357354307251
$ws.Cells.Item(7, 10).Value = 0.9329357354307249
$ws.Cells.Item(7, 13).Value = 0.6933189999999999
$ws.Cells.Item(7, 14).Value = 2.079957
$ws.Cells.Item(7, 15).Value = 0.07063193971031816
$ws.Cells.Item(7, 16).Value = 0.07063193971031817
$ws.Cells.Item(7, 17).Value = 17.25731979961366
$ws.Cells.Item(7, 18).Value = 155.315878196523
$ws.Cells.Item(7, 19).Value = 0.0658950606185443
$ws.Cells.Item(7, 20).Value = 0.0658950606185443

$ws.Cells.Item(8, 7).Value = 1.490294
$ws.Cells.Item(8, 8).Value = 4.470882
$ws.Cells.Item(8, 9).Value = 0.05585774980704767
$ws.Cells.Item(8, 10).Value = 0.05585774980704766
$ws.Cells.Item(8, 13).Value = 0.110552
$ws.Cells.Item(8, 14).Value = 0.331656
$ws.Cells.Item(8, 15).Value = 0.01126249561724847
$ws.Cells.Item(8, 16).Value = 0.01126249561724847
$ws.Cells.Item(8, 17).Value = 0.164754982288
$ws.Cells.Item(8, 18).Value = 1.482794840592
$ws.Cells.Item(8, 19).Value = 0.0006290976623912361
$ws.Cells.Item(8, 20).Value = 0.0006290976623912361

$ws.Cells.Item(9, 7).Value = 1.490294
$ws.Cells.Item(9, 8).Value = 4.470882
$ws.Cells.Item(9, 9).Value = 0.05585774980704767
$ws.Cells.Item(9, 10).Value = 0.05585774980704766
$ws.Cells.Item(9, 15).Value = 0.9181055646724333
$ws.Cells.Item(9, 16).Value = 0.9181055646724334
$ws.Cells.Item(9, 17).Value = 13.43063484210933
$ws.Cells.Item(9, 18).Value = 120.875713578984
$ws.Cells.Item(9, 19).Value = 0.051283310927931
$ws.Cells.Item(9, 20).Value = 0.051283310927931

$ws.Cells.Item(10, 7).Value = 1.490294
$ws.Cells.Item(10, 8).Value = 4.470882
$ws.Cells.Item(10, 9).Value = 0.05585774980704767
$ws.Cells.Item(10, 10).Value = 0.05585774980704766
$ws.Cells.Item(10, 13).Value = 0.6933189999999999
$ws.Cells.Item(10, 14).Value = 2.079957
$ws.Cells.Item(10, 15).Value = 0.07063193971031816
$ws.Cells.Item(10, 16).Value = 0.07063193971031817
$ws.Cells.Item(10, 17).Value = 1.033249145786
$ws.Cells.Item(10, 18).Value = 9.299242312073998
$ws.Cells.Item(10, 19).Value = 0.003945341216725426
$ws.Cells.Item(10, 20).Value = 0.003945341216725427
